# Auto-generated edit script: updates Golem Profits (Leve profit) data
# across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 235.33333
$ws.Range("I2").Value = 300
$ws.Range("J2").Value = 203
$ws.Range("K2").Value = 300
$ws.Range("L2").Value = 203
$ws.Range("M2").Value = -187
$ws.Range("N2").Value = -429

$ws.Range("H29").Value = 791.6667
$ws.Range("J29").Value = 1000
$ws.Range("L29").Value = 3000
$ws.Range("N29").Value = -3562

$ws.Range("H33").Value = 362.5625
$ws.Range("I33").Value = 320.2
$ws.Range("K33").Value = 320.2
$ws.Range("M33").Value = -91.19999999999999

$ws.Range("H95").Value = 40474.75
$ws.Range("J95").Value = 40474.75
$ws.Range("L95").Value = 40474.75
$ws.Range("N95").Value = -45966.75

$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

$ws.Range("H111").Value = 674
$ws.Range("I111").Value = 467.5
$ws.Range("J111").Value = 1500
$ws.Range("K111").Value = 1402.5
$ws.Range("L111").Value = 4500
$ws.Range("M111").Value = 1664.5
$ws.Range("N111").Value = -10634

$ws.Range("H116").Value = 500
$ws.Range("J116").Value = 500
$ws.Range("L116").Value = 500
$ws.Range("N116").Value = -7384

$ws.Range("H140").Value = 52000
$ws.Range("I140").Value = 40000
$ws.Range("K140").Value = 40000
$ws.Range("M140").Value = -34820

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H50").Value = 18246.75
$ws.Range("J50").Value = 34892.5
$ws.Range("L50").Value = 34892.5
$ws.Range("N50").Value = -36320.5

$ws.Range("H61").Value = 597.8
$ws.Range("I61").Value = 597.8
$ws.Range("K61").Value = 597.8
$ws.Range("M61").Value = -385.8

$ws.Range("H92").Value = 59999.75
$ws.Range("J92").Value = 59999.75
$ws.Range("L92").Value = 59999.75
$ws.Range("N92").Value = -64991.75

$ws.Range("H97").Value = 1929
$ws.Range("I97").Value = 1364
$ws.Range("J97").Value = 2720
$ws.Range("K97").Value = 1364
$ws.Range("L97").Value = 2720
$ws.Range("M97").Value = -868
$ws.Range("N97").Value = -3712

$ws.Range("H122").Value = 2500
$ws.Range("I122").Value = 2500
$ws.Range("K122").Value = 7500
$ws.Range("M122").Value = -5050

$ws.Range("H136").Value = 597.8
$ws.Range("I136").Value = 597.8
$ws.Range("K136").Value = 1793.4
$ws.Range("M136").Value = 756.6000000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 140000
$ws.Range("J92").Value = 140000
$ws.Range("L92").Value = 140000
$ws.Range("N92").Value = -144992

$ws.Range("H105").Value = 3953.3333
$ws.Range("I105").Value = 3953.3333
$ws.Range("K105").Value = 3953.3333
$ws.Range("M105").Value = -2206.3333

$ws.Range("H134").Value = 507.5
$ws.Range("I134").Value = 507.5
$ws.Range("K134").Value = 1522.5
$ws.Range("M134").Value = 1012.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 548.6087
$ws.Range("I22").Value = 661.13336
$ws.Range("J22").Value = 337.625
$ws.Range("K22").Value = 661.13336
$ws.Range("L22").Value = 337.625
$ws.Range("M22").Value = -311.13336
$ws.Range("N22").Value = -1037.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J80").Value = 1500
$ws.Range("L80").Value = 4500
$ws.Range("N80").Value = -6372

$ws.Range("J83").Value = 1500
$ws.Range("L83").Value = 13500
$ws.Range("N83").Value = -22860

$ws.Range("H131").Value = 5000
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 5000
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 15000
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -25080

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()

$ws.Range("H95").Value = 35000
$ws.Range("J95").Value = 35000
$ws.Range("L95").Value = 35000
$ws.Range("N95").Value = -40492

$ws.Range("H107").Value = 837.9
$ws.Range("I107").Value = 447
$ws.Range("J107").Value = 1750
$ws.Range("K107").Value = 447
$ws.Range("L107").Value = 1750
$ws.Range("M107").Value = 1473
$ws.Range("N107").Value = -5590

$ws.Range("H113").Value = 821.4
$ws.Range("I113").Value = 717.1429000000001
$ws.Range("K113").Value = 717.1429000000001
$ws.Range("M113").Value = 1452.8571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 3005
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 3005
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 3005
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -3229

$ws.Range("H15").Value = 3005
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 3005
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 3005
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value = -3345

$ws.Range("H46").Value = 258310.88
$ws.Range("J46").Value = 9498.143
$ws.Range("L46").Value = 9498.143
$ws.Range("N46").Value = -9874.143

$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()

$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()

$ws.Range("H93").Value = 2416.375
$ws.Range("I93").Value = 2271.8333
$ws.Range("J93").Value = 2850
$ws.Range("K93").Value = 2271.8333
$ws.Range("L93").Value = 2850
$ws.Range("M93").Value = -1023.8333
$ws.Range("N93").Value = -5346

$ws.Range("H100").Value = 1298.375
$ws.Range("I100").Value = 1298.375
$ws.Range("K100").Value = 1298.375
$ws.Range("M100").Value = -757.375

$ws.Range("H122").Value = 2802
$ws.Range("I122").Value = 2802
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8406
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5956
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 2342
$ws.Range("I132").Value = 1701.3334
$ws.Range("J132").Value = 3303
$ws.Range("K132").Value = 5104.0002
$ws.Range("L132").Value = 9909
$ws.Range("M132").Value = -2574.0002
$ws.Range("N132").Value = -14969

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H117").Value = 49444
$ws.Range("J117").Value = 49444
$ws.Range("L117").Value = 49444
$ws.Range("N117").Value = -58622

$ws.Range("H132").Value = 899.4
$ws.Range("I132").Value = 899.4
$ws.Range("K132").Value = 2698.2
$ws.Range("M132").Value = -168.1999999999998

$ws.Range("H136").Value = 650.8182
$ws.Range("I136").Value = 435.9
$ws.Range("K136").Value = 1307.7
$ws.Range("M136").Value = 1242.3

$ws.Range("H140").Value = 65000
$ws.Range("J140").Value = 65000
$ws.Range("L140").Value = 65000
$ws.Range("N140").Value = -75360
